$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for the MAD-based outlier detection columns
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy header formatting (bold, centered, bordered) from an existing header cell
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# Fill the new columns with boolean FALSE values for every data row (2-25)
$ws.Range("F2:H25").Value = $false
